# Update CCDate timestamps on test-result sheets (Autopay Test cases done)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NoModifyAmountCC")
$ws.Range("B2").Value = "Sat Oct 05 00:25:44 IST 2024"

$ws = $wb.Worksheets.Item("NoModifyAmountCorp")
$ws.Range("B2").Value = "Sat Oct 05 00:26:49 IST 2024"

$ws = $wb.Worksheets.Item("NoModifyAmountPC")
$ws.Range("B2").Value = "Sat Oct 05 00:27:54 IST 2024"

$ws = $wb.Worksheets.Item("NoModifyAmountPS")
$ws.Range("B2").Value = "Sat Oct 05 00:29:02 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayCorp")
$ws.Range("B2").Value = "Sat Oct 05 00:30:09 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:31:16 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayCredit")
$ws.Range("B2").Value = "Sat Oct 05 00:32:22 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:33:29 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayPC")
$ws.Range("B2").Value = "Sat Oct 05 00:34:35 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:35:43 IST 2024"

$ws = $wb.Worksheets.Item("OverAndUnderPayPS")
$ws.Range("B2").Value = "Sat Oct 05 00:36:51 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:37:58 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Sat Oct 05 01:12:23 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCorp")
$ws.Range("B2").Value = "Sat Oct 05 00:41:57 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:43:10 IST 2024"
$ws.Range("B4").Value = "Sat Oct 05 00:44:23 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCorpDCF")
$ws.Range("B2").Value = "Sat Oct 05 00:39:05 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:40:34 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCorpSCF")
$ws.Range("B2").Value = "Sat Oct 05 00:45:33 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:47:01 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCreditDCF")
$ws.Range("B2").Value = "Sat Oct 05 00:48:23 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:49:45 IST 2024"

$ws = $wb.Worksheets.Item("PayNowCreditSCF")
$ws.Range("B2").Value = "Sat Oct 05 00:51:14 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:52:35 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPC")
$ws.Range("B2").Value = "Sat Oct 05 00:54:02 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 00:55:13 IST 2024"
$ws.Range("B4").Value = "Sat Oct 05 00:56:17 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPS")
$ws.Range("B2").Value = "Sat Oct 05 01:09:00 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 01:10:10 IST 2024"
$ws.Range("B4").Value = "Sat Oct 05 01:11:15 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckDCF")
$ws.Range("B2").Value = "Sat Oct 05 01:03:15 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 01:04:44 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPersonalCheckSCF")
$ws.Range("B2").Value = "Sat Oct 05 01:00:20 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 01:01:50 IST 2024"

$ws = $wb.Worksheets.Item("PayNowPersonalSavingsSCF")
$ws.Range("B2").Value = "Sat Oct 05 01:06:08 IST 2024"
$ws.Range("B3").Value = "Sat Oct 05 01:07:36 IST 2024"

$ws = $wb.Worksheets.Item("VerifyConfirmPageLabelsCorp")
$ws.Range("B2").Value = "Sat Oct 05 01:25:15 IST 2024"

$ws = $wb.Worksheets.Item("VerifyConfirmPageLabelsCredit")
$ws.Range("B2").Value = "Mon Oct 07 19:40:04 IST 2024"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsCorp")
$ws.Range("B2").Value = "Sat Oct 05 01:50:00 IST 2024"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsCredit")
$ws.Range("B2").Value = "Sat Oct 05 01:48:48 IST 2024"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsPC")
$ws.Range("B2").Value = "Sat Oct 05 01:51:13 IST 2024"

$ws = $wb.Worksheets.Item("VerifyReceiptPageLabelsPS")
$ws.Range("B2").Value = "Sat Oct 05 01:52:36 IST 2024"
